# Applies scheduled Kraken market-data profit recalculation updates
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

# ALC row 17: One for the Road
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1748.3334
$ws.Range("J17").Value = 1760
$ws.Range("L17").Value = 5280
$ws.Range("N17").Value = -5616

# ALC row 28: The Writing Is Not on the Wall
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 925.2222
$ws.Range("I28").Value = 1331
$ws.Range("J28").Value = 418
$ws.Range("K28").Value = 1331
$ws.Range("L28").Value = 418
$ws.Range("M28").Value = -846
$ws.Range("N28").Value = -1388

# ALC row 39: Riches' Brew
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 88.59999999999999
$ws.Range("I39").Value = 48.25
$ws.Range("J39").Value = 250
$ws.Range("K39").Value = 144.75
$ws.Range("L39").Value = 750
$ws.Range("M39").Value = 151.25
$ws.Range("N39").Value = -1342

# ALC row 51: A Bile Business
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()

# ALC row 132: Fast-forwarding Flora
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6500.8335
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470

# ALC row 135: For Tired Minds
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2065.4285
$ws.Range("J135").Value = 2194
$ws.Range("L135").Value = 19746
$ws.Range("N135").Value = -24816

# ALC row 137: Cutting Edge of Culinary Quality
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1839.3529
$ws.Range("I137").Value = 1779.3334
$ws.Range("K137").Value = 5338.0002
$ws.Range("M137").Value = -2788.0002

# ALC row 138: All-night Crafting
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2983.3215
$ws.Range("I138").Value = 1155.2
$ws.Range("J138").Value = 3998.9443
$ws.Range("K138").Value = 3465.6
$ws.Range("L138").Value = 11996.8329
$ws.Range("M138").Value = 1674.4
$ws.Range("N138").Value = -22276.8329

# ARM row 25: Still Crazy After All These Years
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 3000
$ws.Range("I25").Value = 3000
$ws.Range("K25").Value = 3000
$ws.Range("M25").Value = -2598

# ARM row 61: Dealing with the Tough Stuff
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3981.6667
$ws.Range("I61").Value = 4078
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 4078
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -3866
$ws.Range("N61").Value = -3924

# ARM row 74: As the Bolt Flies
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3639.1667
$ws.Range("I74").Value = 3205.3333
$ws.Range("J74").Value = 4073
$ws.Range("K74").Value = 3205.3333
$ws.Range("L74").Value = 4073
$ws.Range("M74").Value = -2331.3333
$ws.Range("N74").Value = -5821

# ARM row 77: Heavy Metal Banned (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3639.1667
$ws.Range("I77").Value = 3205.3333
$ws.Range("J77").Value = 4073
$ws.Range("K77").Value = 16026.6665
$ws.Range("L77").Value = 20365
$ws.Range("M77").Value = -11658.6665
$ws.Range("N77").Value = -29101

# ARM row 102: Smells of Rich Tama-hagane
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1166.3334
$ws.Range("I102").Value = 1166.3334
$ws.Range("K102").Value = 1166.3334
$ws.Range("M102").Value = 455.6666

# ARM row 132: Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3000
$ws.Range("I132").Value = 3000
$ws.Range("K132").Value = 9000
$ws.Range("M132").Value = -6470

# ARM row 136: Metal with Mettle
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3981.6667
$ws.Range("I136").Value = 4078
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 12234
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -9684
$ws.Range("N136").Value = -15600

# BSM row 37: That's Some Fine Grinding
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1598.8334
$ws.Range("I37").Value = 418.6
$ws.Range("K37").Value = 418.6
$ws.Range("M37").Value = -281.6

# BSM row 134: Ruthenium Supremium
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9683.5
$ws.Range("I134").Value = 5759.4
$ws.Range("K134").Value = 17278.2
$ws.Range("M134").Value = -14743.2

# CRP row 13: Compulsory Conjury
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

# CRP row 58: You Do the Heavy Lifting
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3932.111
$ws.Range("I58").Value = 4499.143
$ws.Range("J58").Value = 1947.5
$ws.Range("K58").Value = 4499.143
$ws.Range("L58").Value = 1947.5
$ws.Range("M58").Value = -4296.143
$ws.Range("N58").Value = -2353.5

# CRP row 132: Hull Lotta Damage
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1100.1666
$ws.Range("I132").Value = 925.25
$ws.Range("J132").Value = 1450
$ws.Range("K132").Value = 2775.75
$ws.Range("L132").Value = 4350
$ws.Range("M132").Value = -245.75
$ws.Range("N132").Value = -9410

# CRP row 136: Turali Quality
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3932.111
$ws.Range("I136").Value = 4499.143
$ws.Range("J136").Value = 1947.5
$ws.Range("K136").Value = 13497.429
$ws.Range("L136").Value = 5842.5
$ws.Range("M136").Value = -10947.429
$ws.Range("N136").Value = -10942.5

# CUL row 36: Love's Crumpets Lost
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 1832.3334
$ws.Range("J36").Value = 3500
$ws.Range("L36").Value = 10500
$ws.Range("N36").Value = -10838

# CUL row 113: Can't Eat Just One
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 944
$ws.Range("I113").Value = 859
$ws.Range("K113").Value = 2577
$ws.Range("M113").Value = -407

# GSM row 122: Awarding Academic Excellence
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2491.1052
$ws.Range("I122").Value = 2416.4285
$ws.Range("K122").Value = 7249.2855
$ws.Range("M122").Value = -4799.2855

# LTW row 7: Tan Before the Ban
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7815.5
$ws.Range("I7").Value = 7878.6
$ws.Range("K7").Value = 7878.6
$ws.Range("M7").Value = -7766.6

# LTW row 68: You Could Say It's a Moving Target
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2566.2222
$ws.Range("I68").Value = 2566.2222
$ws.Range("K68").Value = 2566.2222
$ws.Range("M68").Value = -1817.2222

# LTW row 71: They Call It Bloody Mary (L)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2566.2222
$ws.Range("I71").Value = 2566.2222
$ws.Range("K71").Value = 12831.111
$ws.Range("M71").Value = -9087.111000000001

# LTW row 122: Hell on Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7333.3335
$ws.Range("I122").Value = 8000
$ws.Range("K122").Value = 24000
$ws.Range("M122").Value = -21550

# LTW row 126: Battered Books
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7815.5
$ws.Range("I126").Value = 7878.6
$ws.Range("K126").Value = 23635.8
$ws.Range("M126").Value = -21165.8

# WVR row 136: Weaving the Envelope
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2086.2856
$ws.Range("I136").Value = 2340.8
$ws.Range("J136").Value = 1450
$ws.Range("K136").Value = 7022.400000000001
$ws.Range("L136").Value = 4350
$ws.Range("M136").Value = -4472.400000000001
$ws.Range("N136").Value = -9450
